$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.375107526779175
$ws.Range("B1").Value = 2.850160837173462
$ws.Range("C1").Value = 2.579847574234009
$ws.Range("D1").Value = 2.890549421310425
$ws.Range("E1").Value = 2.918912172317505
